$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the summary header fields (totals recalculated for the new dataset) ---
$ws.Range("E11").Value = 1584148
$ws.Range("C13").Value = 9
$ws.Range("F13").Value = 18

# --- Copy the "last row" emphasized formatting (bottom border style) from the
#     current last data row (43) onto what will become the new last data row
#     (38) before we overwrite/delete anything ---
$ws.Range("B43:J43").Copy()
$ws.Range("B38:J38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Overwrite rows 16-38 with the updated account-statement detail rows ---
function Set-RowData($sheet, $r, $tipoDoc, $numDoc, $nombre, $periodo, $valorMora, $salario) {
    $sheet.Range("B$r").Value = $tipoDoc
    $sheet.Range("C$r").Value = $numDoc
    $sheet.Range("D$r").Value = $nombre
    $sheet.Range("E$r").Value = $periodo
    $sheet.Range("F$r").Value = $valorMora
    $sheet.Range("G$r").Value = $salario
}

Set-RowData $ws 16 "CC" "1052960119" "SHEYLA ANDREA PEREZ MIRANDA" "1709" 29509 737717
Set-RowData $ws 17 "CC" "1050039005" "JUAN PABLO BUELVAS LEYVA" "1802" 31249 781242
Set-RowData $ws 18 "CC" "1143401543" "ISAURA FILO ARBOLEDA" "1905" 6625 828116
Set-RowData $ws 19 "CC" "20373392" "LILIANA VANEGAS ORTEGA" "2207" 68658 1980500
Set-RowData $ws 20 "CC" "79427140" "RAMON JOSE ARIZA RIOS" "2208" 88000 2200000
Set-RowData $ws 21 "CC" "1050969488" "ARLEY ALEXANDER MACIAS TRESPALACIOS" "2208" 40000 1000000
Set-RowData $ws 22 "CC" "1052998110" "MARIAN ALEJANDRA BENAVIDEZ ACOSTA" "2208" 40000 1000000
Set-RowData $ws 23 "CC" "79427140" "RAMON JOSE ARIZA RIOS" "2209" 88000 2200000
Set-RowData $ws 24 "CC" "79427140" "RAMON JOSE ARIZA RIOS" "2210" 88000 2200000
Set-RowData $ws 25 "CC" "79427140" "RAMON JOSE ARIZA RIOS" "2211" 88000 2200000
Set-RowData $ws 26 "CC" "79427140" "RAMON JOSE ARIZA RIOS" "2212" 88000 2200000
Set-RowData $ws 27 "CC" "79427140" "RAMON JOSE ARIZA RIOS" "2301" 88000 2200000
Set-RowData $ws 28 "CC" "79427140" "RAMON JOSE ARIZA RIOS" "2302" 88000 2200000
Set-RowData $ws 29 "CC" "73127846" "JOSE LUIS CARRILLO GRAU" "2302" 46400 1160000
Set-RowData $ws 30 "CC" "79427140" "RAMON JOSE ARIZA RIOS" "2303" 88000 2200000
Set-RowData $ws 31 "CC" "1007939014" "SHARON JULIANA GUERRERO RODRIGUEZ" "2303" 43307 1160000
Set-RowData $ws 32 "CC" "79427140" "RAMON JOSE ARIZA RIOS" "2304" 88000 2200000
Set-RowData $ws 33 "CC" "1007939014" "SHARON JULIANA GUERRERO RODRIGUEZ" "2304" 46400 1160000
Set-RowData $ws 34 "CC" "79427140" "RAMON JOSE ARIZA RIOS" "2305" 88000 2200000
Set-RowData $ws 35 "CC" "79427140" "RAMON JOSE ARIZA RIOS" "2306" 88000 2200000
Set-RowData $ws 36 "CC" "79427140" "RAMON JOSE ARIZA RIOS" "2307" 88000 2200000
Set-RowData $ws 37 "CC" "79427140" "RAMON JOSE ARIZA RIOS" "2308" 88000 2200000
Set-RowData $ws 38 "CC" "79427140" "RAMON JOSE ARIZA RIOS" "2309" 88000 2200000

# --- Remove the now-obsolete trailing detail rows (39-43); this shifts the
#     closing signature block (old rows 48-49) up to rows 43-44 ---
$ws.Rows("39:43").Delete()
